$d = $word.ActiveDocument

$pairs = @(
    @{old="430÷8="; new="548÷4="},
    @{old="904÷8="; new="851÷7="},
    @{old="527÷2="; new="792÷5="},
    @{old="794÷8="; new="882÷2="},
    @{old="383÷8="; new="230÷6="},
    @{old="846÷2="; new="184÷3="},
    @{old="717÷2="; new="593÷7="},
    @{old="598÷9="; new="190÷8="},
    @{old="514÷8="; new="808÷8="},
    @{old="476÷2="; new="425÷4="},
    @{old="564÷2="; new="622÷7="},
    @{old="625÷3="; new="862÷6="},
    @{old="762÷7="; new="278÷4="},
    @{old="390÷8="; new="429÷8="},
    @{old="696÷4="; new="720÷4="},
    @{old="814÷5="; new="410÷6="},
    @{old="333÷7="; new="183÷4="},
    @{old="192÷2="; new="927÷4="},
    @{old="204÷4="; new="467÷2="},
    @{old="856÷7="; new="668÷3="},
    @{old="542÷4="; new="656÷2="},
    @{old="813÷3="; new="521÷8="},
    @{old="187÷2="; new="993÷5="},
    @{old="478÷9="; new="992÷7="},
    @{old="449÷3="; new="886÷2="}
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.new, 2)
}
